# Update the "Förändrad" (Changed) date column for rows 2-5 from
# 45204 (2023-10-05) to 45207 (2023-10-08), matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

foreach ($row in 2..5) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45204) {
        $cell.Value = 45207
    }
}
